$d = $word.ActiveDocument

# Helper: locate $old (case-sensitive, whole-match) and swap in $new while
# keeping the paragraph's leading empty run (<w:r/>) intact. A plain
# Find.Execute(..., Replace:=wdReplaceAll) rewrites the matched run with the
# default (no-rPr) formatting, which is identical to a preceding empty run's
# formatting, so the engine merges the two into one run and the empty run
# disappears. Toggling Bold on/off around the text assignment keeps the
# replacement run's formatting state distinguishable at write time, so the
# preceding empty run survives; the Bold toggle cancels itself out and
# leaves no visible formatting behind.
function Replace-Preserving($doc, $old, $new) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Bold = 1
        $rng.Text = $new
        $rng.Bold = 0
    }
    return $found
}

# 1. Main heading title (appears twice: the H1 heading and the bold run near the end).
#    Find.Execute with Replace:=wdReplaceAll (2) over the full document Content replaces every match in one pass.
$d.Content.Find.Execute("Play Chicken Fox Free: Cute Farm-Inspired Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "Play Chicken Fox for Free - Exciting Farm-Inspired Slot Game", 2)

# 2. "What we like" bullet list items
Replace-Preserving $d "Cute farm-inspired graphics" "Cartoon-style graphics capture the countryside feel"
Replace-Preserving $d "Exciting bonus features" "Exciting bonus features including multipliers and free spins"
Replace-Preserving $d "High payout potential" "High payout potential of up to 13,500 times the bet per round"
Replace-Preserving $d "Solid return potential" "Entertaining gameplay and cute, playful design"

# 3. "What we don't like" bullet list items
Replace-Preserving $d "Medium variance may slow down action" "Medium variance may slow down big wins"
Replace-Preserving $d "Minimalist graphics may not appeal to some players" "Graphics are minimalist"

# 4. Closing summary italic run
$d.Content.Find.Execute("Check out our review of Chicken Fox, a cute and playful online slot game with exciting bonus features. Play for free and enjoy the farm-inspired graphics.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Chicken Fox, a cute and playful online slot game with rewarding features. Play for free and enjoy the farm-inspiration theme.", 2)
